$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = 97389
$ws.Range("B2").Value = 249524
$ws.Range("B3").Value = 97436
$ws.Range("B4").Value = 249593

$excel.Calculate()
